# Add four new "2A motif" rows (GSG-T2A/P2A/E2A/F2A) with their linker
# sequences to the Sequences sheet, rows 76-79, columns A:C and E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row data -----------------------------------------------------
$motifs = @(
    @{ Name = "GSG-T2A"; Seq = "GGAAGCGGAgagggcagaggcagtctgctgacatgcggtgacgtggaagagaatcccggccct" },
    @{ Name = "GSG-P2A"; Seq = "GGAAGCGGAgccaccaacttctccctgctgaagcaggccggcgacgtggaggagaaccccggcccc" },
    @{ Name = "GSG-E2A"; Seq = "GGAAGCGGAcagtgtactaattatgctctcttgaaattggctggagatgttgagagcaacccaggtccc" },
    @{ Name = "GSG-F2A"; Seq = "GGAAGCGGAgtgaaacagactttgaattttgaccttctcaagttggcgggagacgtggagtccaaccctggacct" }
)

$startRow = 76

# Shared-string table order must come out as: all 4 sequences first
# (column C), then all 4 short names (column B) - matching how the
# workbook was actually authored - so fill column-by-column rather than
# row-by-row.
for ($i = 0; $i -lt $motifs.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 3).Value = $motifs[$i].Seq
}
for ($i = 0; $i -lt $motifs.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $motifs[$i].Name
}
for ($i = 0; $i -lt $motifs.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "2A motif"
    $ws.Cells.Item($r, 5).Value = 1
}

$endRow = $startRow + $motifs.Count - 1

# --- Formatting for the new rows --------------------------------------
# Columns A:C -> vertical-center alignment, new black Calibri font
$abc = $ws.Range("A" + $startRow + ":C" + $endRow)
$abc.Font.Name = "Calibri"
$abc.Font.Color = 0
$abc.VerticalAlignment = -4108

# Column E -> same font, right + vertical-center alignment
$colE = $ws.Range("E" + $startRow + ":E" + $endRow)
$colE.Font.Name = "Calibri"
$colE.Font.Color = 0
$colE.VerticalAlignment = -4108
$colE.HorizontalAlignment = -4152

# --- View state: scroll position / active selection --------------------
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$excel.Goto($ws.Range("C81"), $true)

Write-Output "Added rows 76-79 with 2A motif sequences"
